$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the row-label strings on the "BAU Emissions" sheet.
#    Every label of the form "... : NoSettings" becomes "... : test".
#    These labels live in column A, rows 4-28 and 31-280.
# ---------------------------------------------------------------------------
$wsBau = $wb.Worksheets.Item("BAU Emissions")

function Update-NoSettingsLabels($startRow, $endRow) {
    for ($r = $startRow; $r -le $endRow; $r++) {
        $cell = $wsBau.Cells.Item($r, 1)
        $old = $cell.Value2
        if ($old -ne $null) {
            $oldText = $old.ToString()
            if ($oldText.EndsWith(" : NoSettings")) {
                $newText = $oldText.Replace(" : NoSettings", " : test")
                $cell.Value = $newText
            }
        }
    }
}

Update-NoSettingsLabels 4 28
Update-NoSettingsLabels 31 280

# ---------------------------------------------------------------------------
# 2. Correct the data values for row 94 (columns M through AE).
# ---------------------------------------------------------------------------
$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300
$wsBau.Range("Q94").Value = 5005380
$wsBau.Range("R94").Value = 5005380
$wsBau.Range("S94").Value = 5005380
$wsBau.Range("T94").Value = 5005380
$wsBau.Range("U94").Value = 5005380
$wsBau.Range("V94").Value = 5005380
$wsBau.Range("W94").Value = 5005380
$wsBau.Range("X94").Value = 5005380
$wsBau.Range("Y94").Value = 5005380
$wsBau.Range("Z94").Value = 5005380
$wsBau.Range("AA94").Value = 5005380
$wsBau.Range("AB94").Value = 5005380
$wsBau.Range("AC94").Value = 5005380
$wsBau.Range("AD94").Value = 5005380
$wsBau.Range("AE94").Value = 5005380

# ---------------------------------------------------------------------------
# 3. Update the "About" sheet date cell (C1).
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# ---------------------------------------------------------------------------
# 4. Refresh the view/selection state of the sheets that changed:
#    - BAU Emissions: selection becomes A30:AE280
#    - About: becomes the active (selected) sheet/tab
#    Activating "About" last makes it the workbook's active tab, which also
#    removes the previously-selected tab flag from "Current and Planned
#    Capacity".
# ---------------------------------------------------------------------------
$wsBau.Activate()
$wsBau.Range("A30:AE280").Select()

$wsAbout.Activate()
